# Update column G ("K") values for rows 2-24 on Sheet1, per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 3
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 3
    12 = 1
    13 = 2
    14 = 0
    15 = 3
    16 = 0
    17 = 0
    18 = 0
    19 = 2
    20 = 2
    21 = 2
    23 = 2
    24 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
